# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above the current row 2 (shifts existing data rows 2-21 down to 3-22)
$ws.Range("A2:C2").Insert()

# Remove the formatting the new row inherited from the row above so it matches
# the plain (unstyled) data rows
$ws.Range("A2:C2").ClearFormats()

# 2. Populate the newly inserted row 2 with its data
$row2 = New-Object 'object[,]' 1,3
$row2[0,0] = -0.4303629223038048
$row2[0,1] = -0.01591238554787755
$row2[0,2] = -0.3245061741155765
$ws.Range("A2:C2").Value = $row2

# 3. Append 9 new data rows after the existing data (now ending at row 22), i.e. rows 23-31
$tailRows = @(
    @(5.814107340924815, -1.352176813518396, -1.502328406361931),
    @(-0.08294376405329906, -0.04769127909872317, 0.643018007278366),
    @(-4.552947707536836, -0.8177924519326552, 2.338970492867865),
    @(-3.988243347456473, 3.882105143631221, -0.1545383128799389),
    @(-1.513886836646995, 12.12395986789412, -6.341941402739815),
    @(4.529875355608262, 0.5813608329838935, -4.973177166546146),
    @(5.40816806745133, -3.943533415554054, -4.557831178192283),
    @(0.5592678154216362, 2.892510315951151, -1.382000362171883),
    @(-7.637950965837174, -7.75840919158016, -1.187092877235683)
)

$startRow = 23
for ($i = 0; $i -lt $tailRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $tailRows[$i]

    $arr = New-Object 'object[,]' 1,3
    $arr[0,0] = $rowData[0]
    $arr[0,1] = $rowData[1]
    $arr[0,2] = $rowData[2]

    $ws.Range("A" + $rowIndex + ":C" + $rowIndex).Value = $arr
}
